# Add a new bookmark row (row 14) to the BookmarksTestData sheet:
#   ID=16, Title="Oracle JavaFX CSS reference",
#   URL="https://docs.oracle.com/javafx/2/api/javafx/scene/doc-files/cssref.html",
#   Description="A good reference for the JavaFX CSS - including colours, fonts etc.",
#   Category="Coding" (reusing the existing "Coding" category text used by other rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A ("ID") holds numeric-looking IDs stored as text in every other
# row (e.g. "1", "2", "4" ...). Force this cell to text first so "16" is
# written as a string rather than being auto-coerced to a number, then
# drop the temporary number-format again so the cell keeps the sheet's
# default (unstyled) look, matching the rest of the column.
$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = "16"
$ws.Range("A14").ClearFormats()

$ws.Range("B14").Value = "Oracle JavaFX CSS reference"
$ws.Range("C14").Value = "https://docs.oracle.com/javafx/2/api/javafx/scene/doc-files/cssref.html"
$ws.Range("D14").Value = "A good reference for the JavaFX CSS - including colours, fonts etc."
$ws.Range("E14").Value = "Coding"
